$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '56.818.61'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.07%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.394.44'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.55%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '503.97'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.13%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '132.23'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.76%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.399.04'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.09%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0974'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.05%  '
$ws.Range("E11").Value = '  -1.22%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.322'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.31%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.57'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.86%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.820.25'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.49%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '56.752.57'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.06%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.75'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.36%  '
$ws.Range("E17").Value = '  +2.07%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.410.95'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.65%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.22'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.55%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.04'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.37%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '309.45'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.46%  '
$ws.Range("E22").Value = '  +1.73%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.84'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.15%  '
$ws.Range("E24").Value = '  +0.37%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '66.54'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.65%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.998'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.43%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.375'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.88%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.152'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.41%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.43'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.65%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '175.22'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.52%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0₃0724'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.58%  '
$ws.Range("E32").Value = '  -0.21%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.12'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.96%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.88'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.35%  '
$ws.Range("E35").Value = '  +0.15%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.996'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.05%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '17.90'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.59%  '
$ws.Range("E38").Value = '  -0.59%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.81'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.20%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '36.84'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.34%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.823'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.65%  '
$ws.Range("E42").Value = '  +0.88%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '132.98'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +5.43%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.36'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.70%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.84'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.34%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.565'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.43%  '
$ws.Range("B47").Value = 'Stellar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0910'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.08%  '
$ws.Range("B48").Value = 'Bittensor'
$ws.Range("C48").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '249.92'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.66%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0488'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.40%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0211'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.18%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '17.06'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +8.58%  '
